$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update row 5 (IAM004 - Facebook login/logout test case) ---
# Description text gets a new middle bullet about checking existing Neon identity,
# and the Jira id list gets an extra OPQA-2333 entry. Row grows taller to fit text.
$ws.Range("B5").Value = "OPQA-350||OPQA-2333||OPQA-1930"
$ws.Range("C5").Value = "Verify that existing FB user is able to login and logout successfully||Verify that neon should check whether there are any other existing Neon identity with same email id and social account,if the Neon identity does not exists with STeAM account,after signing into Facebook on Neon as first user.||Verify that user can able to Sign out from Neon after successful login in ""Sign In with Facebook"" in Neon"
$ws.Rows.Item(5).RowHeight = 72

# --- Append two new rows (IAM036, IAM037) after the existing last row (36) ---
# Rows 37/38 do not exist yet, so simply copy formatting from row 36 (last
# existing data row) straight down into them instead of doing an Insert
# (which would shift existing data we don't want shifted).
$ws.Range("A36:E36").Copy($ws.Range("A37:E37"))
$ws.Range("A36:E36").Copy($ws.Range("A38:E38"))

# Row 37 - IAM036
$ws.Range("A37").Value = "IAM036"
$ws.Range("B37").Value = "OPQA-2298"
$ws.Range("C37").Value = "Verify that after successful registration on the NEON landing screen using Facebook, user who already has LinkedIn account with the same emailId are prompted to link their Linked account with the newly created Facebook account"
$ws.Range("D37").Value = "Y"
$ws.Range("E37").Value = "PASS"
$ws.Rows.Item(37).RowHeight = 43.2

# Row 38 - IAM037
$ws.Range("A38").Value = "IAM037"
$ws.Range("B38").Value = "OPQA-1936"
$ws.Range("C38").Value = "Verify that user is able to reset STeAM Password from the Neon landing page.Prerequisites: A Neon identity with a linked STeAM account."
$ws.Range("D38").Value = "Y"
$ws.Range("E38").Value = "PASS"
$ws.Rows.Item(38).RowHeight = 28.8

# --- Update the view so the newly added row 37 is the active selection ---
[void]$ws.Range("E37").Select()
